$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the full row contents (columns B:AB) between two row numbers ---
function Swap-Rows {
    param($ws, $r1, $r2)
    $range1 = $ws.Range($ws.Cells.Item($r1, 2), $ws.Cells.Item($r1, 28))
    $range2 = $ws.Range($ws.Cells.Item($r2, 2), $ws.Cells.Item($r2, 28))
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

Swap-Rows $ws 164 165
Swap-Rows $ws 181 182
Swap-Rows $ws 183 184

# --- Direct odds updates for rows 304-309 ---
$ws.Range("M304").Value2 = 2.15
$ws.Range("O304").Value2 = 3
$ws.Range("Q304").Value2 = 1.925
$ws.Range("R304").Value2 = 1.925
$ws.Range("T304").Value2 = 1.925
$ws.Range("U304").Value2 = 1.925

$ws.Range("M305").Value2 = 2.55
$ws.Range("S305").Value2 = 3.25
$ws.Range("T305").Value2 = 2
$ws.Range("U305").Value2 = 1.85

$ws.Range("T306").Value2 = 1.8
$ws.Range("U306").Value2 = 2.05

$ws.Range("M307").Value2 = 3.4
$ws.Range("N307").Value2 = 3.6
$ws.Range("O307").Value2 = 2
$ws.Range("P307").Value2 = 0.5
$ws.Range("Q307").Value2 = 1.8
$ws.Range("R307").Value2 = 2.05

$ws.Range("M308").Value2 = 2.6
$ws.Range("O308").Value2 = 2.5
$ws.Range("P308").Value2 = 0
$ws.Range("Q308").Value2 = 1.975
$ws.Range("R308").Value2 = 1.875

$ws.Range("Q309").Value2 = 1.85
$ws.Range("R309").Value2 = 2
$ws.Range("T309").Value2 = 1.85
$ws.Range("U309").Value2 = 2
